$d = $word.ActiveDocument

# 1. Update the ID placeholder text in the first paragraph, and drop the
#    trailing-space run that used to follow it (the two runs shared identical
#    formatting, so Word's Find/Replace coalesces them into a single run).
$d.Content.Find.Execute("**ID__AFFARS_5303_topic_14__ID**", $false, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_SUBPART_5303_7__ID**", 2)
$d.Content.Find.Execute("**ID__AFFARS_SUBPART_5303_7__ID** ", $false, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_SUBPART_5303_7__ID**", 2)

# 2. Tweak the first paragraph's formatting: left indent 120 -> 225 twips
#    (i.e. 6pt -> 11.25pt), and add a thin paragraph border (5pt space) on
#    all four sides.
$p1 = $d.Paragraphs(1)
$p1.Range.ParagraphFormat.LeftIndent = 11.25
$p1.Range.ParagraphFormat.Borders.DistanceFromTop = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromLeft = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromBottom = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromRight = 5
